# Natmi following Dr Hou advice
# Rewrites the LR-pair data block (rows 2-7) with updated statistics
# and appends three new rows (8-10) for the "sCs" sending cluster,
# matching the new NATMI run output for Col4a2-Cd93.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Col4a2"
$ws.Cells.Item(2, 3).Value = "Cd93"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 124.506835
$ws.Cells.Item(2, 8).Value = 373.520505
$ws.Cells.Item(2, 9).Value = 0.3727127886796942
$ws.Cells.Item(2, 10).Value = 0.3727127886796942
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 124.722578
$ws.Cells.Item(2, 14).Value = 374.167734
$ws.Cells.Item(2, 15).Value = 0.9767311432246923
$ws.Cells.Item(2, 16).Value = 0.9767311432246923
$ws.Cells.Item(2, 17).Value = 15528.81343982063
$ws.Cells.Item(2, 18).Value = 139759.3209583857
$ws.Cells.Item(2, 19).Value = 0.3640401881815808
$ws.Cells.Item(2, 20).Value = 0.3640401881815808

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Col4a2"
$ws.Cells.Item(3, 3).Value = "Cd93"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 124.506835
$ws.Cells.Item(3, 8).Value = 373.520505
$ws.Cells.Item(3, 9).Value = 0.3727127886796942
$ws.Cells.Item(3, 10).Value = 0.3727127886796942
$ws.Cells.Item(3, 11).Value = 1
$ws.Cells.Item(3, 12).Value = 0.3333333333333333
$ws.Cells.Item(3, 13).Value = 0.134712
$ws.Cells.Item(3, 14).Value = 0.404136
$ws.Cells.Item(3, 15).Value = 0.001054960600366076
$ws.Cells.Item(3, 16).Value = 0.001054960600366076
$ws.Cells.Item(3, 17).Value = 16.77256475652
$ws.Cells.Item(3, 18).Value = 150.95308280868
$ws.Cells.Item(3, 19).Value = 0.0003931973073096446
$ws.Cells.Item(3, 20).Value = 0.0003931973073096446

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Col4a2"
$ws.Cells.Item(4, 3).Value = "Cd93"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 124.506835
$ws.Cells.Item(4, 8).Value = 373.520505
$ws.Cells.Item(4, 9).Value = 0.3727127886796942
$ws.Cells.Item(4, 10).Value = 0.3727127886796942
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 2.836578333333333
$ws.Cells.Item(4, 14).Value = 8.509735
$ws.Cells.Item(4, 15).Value = 0.02221389617494163
$ws.Cells.Item(4, 16).Value = 0.02221389617494163
$ws.Cells.Item(4, 17).Value = 353.1733905129083
$ws.Cells.Item(4, 18).Value = 3178.560514616175
$ws.Cells.Item(4, 19).Value = 0.008279403190803685
$ws.Cells.Item(4, 20).Value = 0.008279403190803685

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Col4a2"
$ws.Cells.Item(5, 3).Value = "Cd93"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 144.226517
$ws.Cells.Item(5, 8).Value = 432.679551
$ws.Cells.Item(5, 9).Value = 0.4317439066909806
$ws.Cells.Item(5, 10).Value = 0.4317439066909806
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 124.722578
$ws.Cells.Item(5, 14).Value = 374.167734
$ws.Cells.Item(5, 15).Value = 0.9767311432246923
$ws.Cells.Item(5, 16).Value = 0.9767311432246923
$ws.Cells.Item(5, 17).Value = 17988.30301620083
$ws.Cells.Item(5, 18).Value = 161894.7271458074
$ws.Cells.Item(5, 19).Value = 0.4216977195625763
$ws.Cells.Item(5, 20).Value = 0.4216977195625763

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Col4a2"
$ws.Cells.Item(6, 3).Value = "Cd93"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 144.226517
$ws.Cells.Item(6, 8).Value = 432.679551
$ws.Cells.Item(6, 9).Value = 0.4317439066909806
$ws.Cells.Item(6, 10).Value = 0.4317439066909806
$ws.Cells.Item(6, 11).Value = 1
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.134712
$ws.Cells.Item(6, 14).Value = 0.404136
$ws.Cells.Item(6, 15).Value = 0.001054960600366076
$ws.Cells.Item(6, 16).Value = 0.001054960600366076
$ws.Cells.Item(6, 17).Value = 19.429042558104
$ws.Cells.Item(6, 18).Value = 174.861383022936
$ws.Cells.Item(6, 19).Value = 0.000455472811007112
$ws.Cells.Item(6, 20).Value = 0.000455472811007112

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Col4a2"
$ws.Cells.Item(7, 3).Value = "Cd93"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 144.226517
$ws.Cells.Item(7, 8).Value = 432.679551
$ws.Cells.Item(7, 9).Value = 0.4317439066909806
$ws.Cells.Item(7, 10).Value = 0.4317439066909806
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.836578333333333
$ws.Cells.Item(7, 14).Value = 8.509735
$ws.Cells.Item(7, 15).Value = 0.02221389617494163
$ws.Cells.Item(7, 16).Value = 0.02221389617494163
$ws.Cells.Item(7, 17).Value = 409.1098132143317
$ws.Cells.Item(7, 18).Value = 3681.988318928985
$ws.Cells.Item(7, 19).Value = 0.009590714317397128
$ws.Cells.Item(7, 20).Value = 0.009590714317397128

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Col4a2"
$ws.Cells.Item(8, 3).Value = "Cd93"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 65.32235733333333
$ws.Cells.Item(8, 8).Value = 195.967072
$ws.Cells.Item(8, 9).Value = 0.1955433046293252
$ws.Cells.Item(8, 10).Value = 0.1955433046293253
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 124.722578
$ws.Cells.Item(8, 14).Value = 374.167734
$ws.Cells.Item(8, 15).Value = 0.9767311432246923
$ws.Cells.Item(8, 16).Value = 0.9767311432246923
$ws.Cells.Item(8, 17).Value = 8147.172807650538
$ws.Cells.Item(8, 18).Value = 73324.55526885485
$ws.Cells.Item(8, 19).Value = 0.1909932354805351
$ws.Cells.Item(8, 20).Value = 0.1909932354805351

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Col4a2"
$ws.Cells.Item(9, 3).Value = "Cd93"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 65.32235733333333
$ws.Cells.Item(9, 8).Value = 195.967072
$ws.Cells.Item(9, 9).Value = 0.1955433046293252
$ws.Cells.Item(9, 10).Value = 0.1955433046293253
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.134712
$ws.Cells.Item(9, 14).Value = 0.404136
$ws.Cells.Item(9, 15).Value = 0.001054960600366076
$ws.Cells.Item(9, 16).Value = 0.001054960600366076
$ws.Cells.Item(9, 17).Value = 8.799705401088
$ws.Cells.Item(9, 18).Value = 79.197348609792
$ws.Cells.Item(9, 19).Value = 0.0002062904820493194
$ws.Cells.Item(9, 20).Value = 0.0002062904820493195

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Col4a2"
$ws.Cells.Item(10, 3).Value = "Cd93"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 65.32235733333333
$ws.Cells.Item(10, 8).Value = 195.967072
$ws.Cells.Item(10, 9).Value = 0.1955433046293252
$ws.Cells.Item(10, 10).Value = 0.1955433046293253
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 2.836578333333333
$ws.Cells.Item(10, 14).Value = 8.509735
$ws.Cells.Item(10, 15).Value = 0.02221389617494163
$ws.Cells.Item(10, 16).Value = 0.02221389617494163
$ws.Cells.Item(10, 17).Value = 185.2919834939911
$ws.Cells.Item(10, 18).Value = 1667.62785144592
$ws.Cells.Item(10, 19).Value = 0.004343778666740813
$ws.Cells.Item(10, 20).Value = 0.004343778666740814

